# Update the "取得日時" (retrieved datetime) column for all data rows on the
# "ランサーズ" sheet from 2025-12-08 06:31:34 to 2025-12-08 06:40:10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-08 06:40:10"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
